$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24 (this shifts existing rows 24-111 down to 25-112,
# matching the target diff where each subsequent row picks up the data that used to
# be one row above it, and the old last row (111) becomes the new last row (112)).
$ws.Rows("24:24").Insert()

# Populate the newly inserted row 24 with the new record's data.
$ws.Range("A24").Value = 4
$ws.Range("B24").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C24").Value = "Los Lagos"
$ws.Range("D24").Value = 44608
$ws.Range("E24").Value = 10
$ws.Range("F24").Value = 100112022
$ws.Range("G24").Value = "Arveja Verde"
$ws.Range("H24").Value = "Sin especificar"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 25
$ws.Range("K24").Value = 20000
$ws.Range("L24").Value = 20000
$ws.Range("M24").Value = 20000
$ws.Range("N24").Value = "$/saco 25 kilos"
$ws.Range("O24").Value = "Región de La Araucanía"
$ws.Range("P24").Value = 800
$ws.Range("Q24").Value = 25
$ws.Range("R24").Value = "Hortaliza"
